$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: show the formula text behind column D's DATE() result ---
$ws.Range("F1").Value = "Formula Text"
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").Font.Bold = $true

$ws.Range("F2").Formula = "=FORMULATEXT(D2)"
$ws.Range("F3:F9").Formula = "=FORMULATEXT(D3)"

# --- Column E ("Return") no longer needs the extra number-format style; just center it ---
$ws.Range("E2:E9").HorizontalAlignment = -4108

# --- Column widths ---
$ws.Columns.Item(5).ColumnWidth = 20.166666666666668
$ws.Columns.Item(6).ColumnWidth = 23.307291666666668

# --- View tweaks ---
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("D9").Select()
